$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row to the table (expands table range A1:D3 -> A1:D4)
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Row 2 - update existing task/owner, keep status "Pronto"
$ws.Range("A2").Value = "03/04/2024 - 10/04/2024"
$ws.Range("B2").Value = "Desenvolvimento das Introduções das Matérias"
$ws.Range("C2").Value = "Grupo"
$ws.Range("D2").Value = "Pronto"

# Row 3 - update existing task/owner, status "Pronto"
$ws.Range("A3").Value = "03/04/2024 - 10/04/2024"
$ws.Range("B3").Value = "Inicio do Manual do Usuário"
$ws.Range("C3").Value = "Gabriel"
$ws.Range("D3").Value = "Pronto"

# Row 4 - new row
$ws.Range("A4").Value = "03/04/2024 - 10/04/2024"
$ws.Range("B4").Value = "Realização dos requisitos do PIM"
$ws.Range("C4").Value = "Grupo"
$ws.Range("D4").Value = "Pronto"

# Data rows now use horizontal-center / (default) vertical alignment,
# instead of horizontal+vertical center
$ws.Range("A2:D4").HorizontalAlignment = -4108
$ws.Range("A2:D4").VerticalAlignment = -4107

# Selection artifact from editing session
$ws.Range("F10").Select() | Out-Null
